$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# --- Insert 12 new rows to expand the old 2-row draft-marks block (rows 16-17)
# into the new 14-row detailed draft-marks table (rows 16-29). Everything
# below (old "Остойчивость" section, previously starting at row 18) shifts
# down by 12 rows (now starts at row 30); formulas/merged cells auto-adjust.
$ws.Rows("18:29").Insert()

# --- New draft-marks rows content ------------------------------------------------
$data = @(
    @(16, "80", "Осадка на кормовых марках ЛБ", "м", "-"),
    @(17, "81", "Осадка на кормовых марках осредненная", "м", "-"),
    @(18, "82", "Осадка на промежуточных кормовых марках ПрБ", "м", 8.6509999999999998),
    @(19, "83", "Осадка на промежуточных кормовых марках ЛБ", "м", "-"),
    @(20, "84", "Осадка на промежуточных кормовых марках осредненная", "м", "-"),
    @(21, "85", "Осадка на миделевых марках ПрБ", "м", "-"),
    @(22, "86", "Осадка на миделевых марках ЛБ", "м", "-"),
    @(23, "87", "Осадка на миделевых марках осредненная", "м", "-"),
    @(24, "88", "Осадка на промежуточных носовых марках ПрБ", "м", "-"),
    @(25, "89", "Осадка на промежуточных носовых марках ЛБ", "м", "-"),
    @(26, "90", "Осадка на промежуточных носовых марках осредненная", "м", "-"),
    @(27, "91", "Осадка на носовых марках ПрБ", "м", "-"),
    @(28, "92", "Осадка на носовых марках ЛБ", "м", 7.2999999999999998),
    @(29, "93", "Осадка на носовых марках осредненная", "м", "-")
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

# Column A on these rows should be centre-aligned like the surrounding table.
$ws.Range("A16:A29").HorizontalAlignment = -4108
